$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5109
$ws.Range("L3").Value = 5508
$ws.Range("I4").Value = 1846
$ws.Range("J4").Value = 1878
$ws.Range("K4").Value = 1788
$ws.Range("L4").Value = 1344
$ws.Range("L5").Value = 325
$ws.Range("L6").Value = 4595
$ws.Range("I7").Value = 26316
$ws.Range("J7").Value = 29355
$ws.Range("K7").Value = 27580
$ws.Range("L7").Value = 16881

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L5").Value = 9
$ws.Range("L6").Value = 20

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 325
$ws.Range("L3").Value = 382
$ws.Range("L6").Value = 288
$ws.Range("L7").Value = 1116

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 153
$ws.Range("L7").Value = 374

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 214
$ws.Range("L3").Value = 272
$ws.Range("L7").Value = 780

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I4").Value = 43
$ws.Range("L4").Value = 34
$ws.Range("I7").Value = 811
$ws.Range("L7").Value = 642

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 120
$ws.Range("L3").Value = 103
$ws.Range("L7").Value = 329

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 144
$ws.Range("L7").Value = 548
$ws.Range("L8").Value = 1116
$ws.Range("L20").Value = 417
$ws.Range("L25").Value = 100
$ws.Range("L29").Value = 931
$ws.Range("L30").Value = 77
$ws.Range("L31").Value = 167
$ws.Range("L33").Value = 780
$ws.Range("L34").Value = 99
$ws.Range("L36").Value = 216
$ws.Range("I37").Value = 811
$ws.Range("L37").Value = 642
$ws.Range("L38").Value = 20
$ws.Range("L40").Value = 47
$ws.Range("L42").Value = 552
$ws.Range("L44").Value = 117
$ws.Range("L47").Value = 116
$ws.Range("L48").Value = 219
$ws.Range("L49").Value = 86
$ws.Range("L50").Value = 86
$ws.Range("L52").Value = 339
$ws.Range("L53").Value = 189
$ws.Range("L54").Value = 361
$ws.Range("L55").Value = 171
$ws.Range("L57").Value = 61
$ws.Range("L61").Value = 18
$ws.Range("L62").Value = 14
$ws.Range("J63").Value = 229
$ws.Range("K63").Value = 174
$ws.Range("L63").Value = 52
$ws.Range("L65").Value = 329
$ws.Range("L66").Value = 43
$ws.Range("L67").Value = 581
$ws.Range("L79").Value = 444
$ws.Range("L83").Value = 374
$ws.Range("L91").Value = 228
$ws.Range("L94").Value = 210
$ws.Range("L95").Value = 237
$ws.Range("L96").Value = 192
$ws.Range("I101").Value = 26316
$ws.Range("J101").Value = 29355
$ws.Range("K101").Value = 27580
$ws.Range("L101").Value = 16881

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 224
$ws.Range("L7").Value = 581

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 361

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 284
$ws.Range("L3").Value = 350
$ws.Range("L4").Value = 47
$ws.Range("L7").Value = 931

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 45
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 154
$ws.Range("L3").Value = 189
$ws.Range("L7").Value = 552

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 52
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 192

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 102
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 444

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 129
$ws.Range("L6").Value = 112
$ws.Range("L7").Value = 417

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 186
$ws.Range("L6").Value = 130
$ws.Range("L7").Value = 548

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 48
$ws.Range("L3").Value = 49
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 21
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 108
$ws.Range("L3").Value = 110
$ws.Range("L7").Value = 339

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 14
